$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated figures from re-running the data prepare & render pipeline
# with the final data (recomputed "All" and "Saudi Arabia"/"USA" columns).

$ws.Range("B2").Value = 18.5840386333021
$ws.Range("K2").Value = 14.7145599857929
$ws.Range("M2").Value = 18.8176682846185

$ws.Range("B3").Value = 18.2790863468141
$ws.Range("K3").Value = 17.7049808040103
$ws.Range("M3").Value = 17.7140126023854

$ws.Range("B4").Value = 16.7223157348188
$ws.Range("K4").Value = 15.9779234610467
$ws.Range("M4").Value = 16.7126947567179

$ws.Range("B5").Value = 16.2405917484416
$ws.Range("K5").Value = 14.0638042866307
$ws.Range("M5").Value = 15.5556817545784
